# Applies the "Trying to get azimuth correct" edit to the Report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# --- New named ranges used by the extended spherical-triangle calc ---
$wb.Names.Add("_PX",  "=Report!`$D`$20")
$wb.Names.Add("_PZ",  "=Report!`$D`$21")
$wb.Names.Add("_ZX",  "=Report!`$D`$22")
$wb.Names.Add("_Alt", "=Report!`$D`$23")
$wb.Names.Add("_ZPA", "=Report!`$D`$24")

# --- Updated inputs (L, d, t) ---
$ws.Range("D7").Value2 = -20
$ws.Range("D8").Value2 = 22.05
$ws.Range("D9").Value2 = 63

# --- New sample column N for rows 7-9, 11, 12 ---
$ws.Range("N7").Value2 = 36
$ws.Range("N8").Value2 = 28.091666666666665
$ws.Range("N9").Value2 = 68
$ws.Range("N11").Formula = "=32+58/60"
$ws.Range("N12").Value2 = 77.099999999999994

# --- New altitude (h) readable DMS string next to D11 ---
$ws.Range("E11").Formula = '=INT(_h)&"° "&TEXT(ROUND(MOD(_h,1)*60,1),"0.0")'

# --- Replace old row 16 (C14+270) with row 15 (180+C14) ---
$ws.Range("C16").ClearContents()
$ws.Range("C15").Formula = "=180+C14"

# --- New spherical-triangle section (rows 20-24) ---
$ws.Range("C20").Value2 = "_PX"
$ws.Range("D20").Formula = "=90-_d"

$ws.Range("C21").Value2 = "_PZ"
$ws.Range("D21").Formula = "=90-_L"

$ws.Range("C22").Value2 = "_ZX"
$ws.Range("D22").Formula = "=90-_h"

$ws.Range("C23").Value2 = "_Alt"
$ws.Range("D23").Formula = "=90-_ZX"

$ws.Range("C24").Value2 = "_ZPA"
$ws.Range("D24").Formula = "=_t"

# --- Reference link near the top of the sheet ---
$ws.Range("H5").Value2 = "https://astronavigationdemystified.com/calculating-azimuth-and-altitude-at-the-assumed-position-by-spherical-trigonometry/"

# --- New cosine-rule based azimuth computation (rows 26-28) ---
$ws.Range("C26").Formula = "=(COS(_PX*PI()/180)-COS(_ZX*PI()/180)*COS(_PZ*PI()/180))/(SIN(_ZX*PI()/180)*SIN(_PZ*PI()/180))"

$ws.Range("C27").Formula = "=ACOS(C26)*180/PI()"
$ws.Range("D27").Formula = '=INT(C27)&"° "&TEXT(ROUND(MOD(C27,1)*60,1),"0.0")'

$ws.Range("C28").Formula = "=180-C27"
$ws.Range("D28").Formula = '=INT(C28)&"° "&TEXT(ROUND(MOD(C28,1)*60,1),"0.0")'

# --- Sheet view: scroll down a bit and move the active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("J19").Select()

$wb.Save()
